$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset the current used range so we can rebuild it cleanly ---
$ws.Range("A1:I9").Clear() | Out-Null

# --- Cell values (sharedStrings content moves around as part of the rework) ---
$ws.Range("A1").Value = "vis"
$ws.Range("B1").Value = "unknown global"
$ws.Range("C1").Value = "helper"
$ws.Range("D1").Value = "hindiChar"
$ws.Range("F1").Value = "poemChar"
$ws.Range("G1").Value = "poemLine"
$ws.Range("H1").Value = "poem"
$ws.Range("I1").Value = "freeVerse"
$ws.Range("J1").Value = "ghazal"
$ws.Range("K1").Value = "needless"
$ws.Range("A2").Value = "charWidth"
$ws.Range("B2").Value = "mode"
$ws.Range("C2").Value = "previousText"
$ws.Range("E2").Value = "index"
$ws.Range("F2").Value = "inherited from: hindiChar"
$ws.Range("G2").Value = "poemChars[]"
$ws.Range("H2").Value = "poemText"
$ws.Range("I2").Value = "inherited from: poem"
$ws.Range("J2").Value = "inherited from: poem"
$ws.Range("K2").Value = "flagFreeVerse"
$ws.Range("A3").Value = "charHeight"
$ws.Range("B3").Value = "prevBaseCount"
$ws.Range("D3").Value = "mainChar"
$ws.Range("E3").Value = 0
$ws.Range("H3").Value = "poemLines[]"
$ws.Range("J3").Value = "radeef"
$ws.Range("K3").Value = "flagGhazal"
$ws.Range("A4").Value = "paddingLeft"
$ws.Range("B4").Value = "selWord"
$ws.Range("D4").Value = "vowelChar"
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = "poem previousVersion"
$ws.Range("J4").Value = "unkown global"
$ws.Range("A5").Value = "lineSpacing"
$ws.Range("B5").Value = "chars[]"
$ws.Range("D5").Value = "consonantNumber"
$ws.Range("E5").Value = 2
$ws.Range("J5").Value = "radeefArray"
$ws.Range("A6").Value = "flagShowText"
$ws.Range("B6").Value = "maxLen"
$ws.Range("D6").Value = "vowelNumber"
$ws.Range("E6").Value = 3
$ws.Range("H6").Value = "new properties"
$ws.Range("J6").Value = "radeefTruncated"
$ws.Range("A7").Value = "flagLineSpacing"
$ws.Range("B7").Value = "maxLineLen = number of lines"
$ws.Range("H7").Value = "title"
$ws.Range("B8").Value = "compositeLinesMarkingA"
$ws.Range("H8").Value = "editDate"
$ws.Range("H9").Value = "poem versions[]"

# --- Header row (row 1): bold, 9pt ---
$ws.Rows.Item(1).Font.Bold = $true
$ws.Rows.Item(1).Font.Size = 9
$ws.Rows.Item(1).Font.ColorIndex = -4105
$ws.Range("A1:C1,G1:K1").Font.Color = 255
$ws.Range("D1:F1").Font.Color = 255
$ws.Range("D1:F1").HorizontalAlignment = -4131

# --- Body cells: 9pt, not bold, default color ---
$ws.Range("A2:D9,G2:G9,H2:H9,J2:K9").Font.Size = 9

# --- "inherited from / index" style cells: 9pt, not bold, red ---
$ws.Range("E2:F2,I2:J2,J4,G6,H6").Font.Size = 9
$ws.Range("E2:F2,I2:J2,J4,G6,H6").Font.Color = 255

# --- Merge the hindiChar header cell ---
$ws.Range("D1:E1").Merge() | Out-Null

# --- Column widths (approximate bestFit look) ---
$ws.Columns.Item(1).ColumnWidth = 12.6
$ws.Columns.Item(2).ColumnWidth = 23.9
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 14.6
$ws.Columns.Item(5).ColumnWidth = 4.5
$ws.Columns.Item(6).ColumnWidth = 20.2
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 18
$ws.Columns.Item(9).ColumnWidth = 17
$ws.Columns.Item(10).ColumnWidth = 17
$ws.Columns.Item(11).ColumnWidth = 11

# --- Selection ---
$ws.Range("G2").Select() | Out-Null
